# Append 8 new GSW box-score rows (A148:AD155) for the 2025-03-30 .. 2025-04-04
# games, per the updated gsw_box_score_team_stats data.
#
# Column layout (row 1 header): A=#, B=TEAM, C=OPP, D=STATUS, E=DATE, F=MIN,
# G=FGM, H=FGA, I=FG%, J=3PM, K=3PA, L=3P%, M=FTM, N=FTA, O=FT%, P=OREB,
# Q=DREB, R=TREB, S=AST, T=STL, U=BLK, V=TOV, W=PF, X=PTS, Y=+/-, Z=Q1,
# AA=Q2, AB=Q3, AC=Q4, AD=W/L.
#
# Existing rows store E (DATE, e.g. "2025-03-30") and F (MIN, e.g. "240:00")
# as literal text, not as Excel dates/times, so those two columns are forced
# to Text format before assignment to stop Excel's autodetection turning
# them into date/time serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    # row, A#,  B,     C,     D,      E,            F,        G,  H,  I,      J,  K,  L,      M,  N,  O,      P,  Q,  R,  S,  T,  U, V,  W,  X,   Y,   Z,  AA, AB, AC, AD
    @(148, 146, "GSW", "SAS", "away", "2025-03-30", "240:00", 56, 97, 0.577, 21, 44, 0.477, 15, 19, 0.789, 10, 34, 44, 42, 17, 3, 12, 15, 148,  42,  44, 24, 43, 37, "W"),
    @(149, 147, "SAS", "GSW", "home", "2025-03-30", "240:00", 40, 84, 0.476, 16, 42, 0.381, 10, 16, 0.625,  7, 27, 34, 30,  4, 3, 21, 20, 106, -42,  27, 17, 29, 33, "L"),
    @(150, 148, "GSW", "MEM", "away", "2025-04-01", "240:00", 42, 99, 0.424, 22, 56, 0.393, 28, 28, 1,     15, 34, 49, 29, 11, 3, 11, 22, 134,   9,  45, 29, 29, 31, "W"),
    @(151, 149, "MEM", "GSW", "home", "2025-04-01", "240:00", 44, 92, 0.478, 18, 44, 0.409, 19, 22, 0.864,  8, 35, 43, 22,  7, 6, 16, 23, 125,  -9,  32, 34, 37, 22, "L"),
    @(152, 150, "GSW", "LAL", "away", "2025-04-03", "240:00", 39, 82, 0.476, 19, 42, 0.452, 26, 31, 0.839, 14, 32, 46, 30,  5, 2, 13, 21, 123,   7,  26, 34, 28, 35, "W"),
    @(153, 151, "LAL", "GSW", "home", "2025-04-03", "240:00", 39, 86, 0.453, 18, 47, 0.383, 20, 23, 0.87,  14, 25, 39, 23,  5, 3,  8, 23, 116,  -7,  22, 25, 30, 39, "L"),
    @(154, 152, "DEN", "GSW", "away", "2025-04-04", "240:00", 40, 74, 0.541, 12, 31, 0.387, 12, 13, 0.923,  7, 31, 38, 28,  8, 4, 25, 18, 104, -14,  44, 16, 24, 20, "L"),
    @(155, 153, "GSW", "DEN", "home", "2025-04-04", "240:00", 43, 88, 0.489, 16, 38, 0.421, 16, 18, 0.889, 11, 25, 36, 25, 14, 2, 13, 15, 118,  14,  34, 32, 30, 22, "W")
)

$numericCols = @("G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($r in $newRows) {
    $row = $r[0]

    # A: row index (number) - styled like the rest of column A below.
    $ws.Range("A$row").Value = $r[1]

    # B/C/D: team, opponent, home-away status - plain text.
    $ws.Range("B$row").Value = $r[2]
    $ws.Range("C$row").Value = $r[3]
    $ws.Range("D$row").Value = $r[4]

    # E/F: date & minutes-played, stored as text like the rest of the sheet.
    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = $r[5]
    $ws.Range("F$row").NumberFormat = "@"
    $ws.Range("F$row").Value = $r[6]

    # G:AC - the 23 numeric box-score stat columns, written in one shot.
    $statVals = $r[7..29]
    $arr = New-Object 'object[,]' 1, $numericCols.Count
    for ($i = 0; $i -lt $numericCols.Count; $i++) { $arr[0, $i] = $statVals[$i] }
    $ws.Range("G$row`:AC$row").Value = $arr

    # AD: W/L result - plain text.
    $ws.Range("AD$row").Value = $r[30]
}

$firstRow = $newRows[0][0]
$lastRow = $newRows[-1][0]

# Column A on the new rows should carry the same bold/bordered/centered
# style as the rest of column A - copy it from the row directly above.
$ws.Range("A147").Copy()
$ws.Range("A$firstRow`:A$lastRow").PasteSpecial(-4122)

# The forced Text number format on E/F above leaves those cells on a new
# style; reset them back to the plain/default style used elsewhere in the
# sheet (values, already text, are untouched by a formats-only paste).
$ws.Range("E147:F147").Copy()
$ws.Range("E$firstRow`:F$lastRow").PasteSpecial(-4122)

$excel.CutCopyMode = $false

Write-Output "Inserted rows $firstRow-$lastRow into $($ws.Name)"
